$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Status / Comentarios for rows 20-22 (DR + cronograma update) ---
$ws.Range("H20").Value = "Concluída"
$ws.Range("I20").Value = "Todas as atividades previstas até o momento foram concluídas."

$ws.Range("H21").Value = "Não houve necessidade"
$ws.Range("I21").Value = "A fase 2 da auditoria do GRE não gerou inconformidades."

$ws.Range("H22").Value = "Em andamento"
$ws.Range("I22").Value = "Documento de requisitos enviado ao cliente para aprovação."

# --- Widen column I slightly to fit the new, longer comments ---
$ws.Columns.Item(9).ColumnWidth = 108.8

# --- Scroll the view over and move the selection down near the bottom ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H24").Select()
